# Applies the "soil_texture" -> "soil_texture_main" / "soil_texture_sub" split
# to the soil_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("soil_data")

# Insert a new column before column F (soil_texture). This shifts the
# existing soil_texture column (and everything after it) one column to
# the right, so it becomes column G.
$ws.Columns.Item(6).Insert()

# Rename headers: new F1 = soil_texture_main, shifted G1 = soil_texture_sub
$ws.Range("F1").Value = "soil_texture_main"
$ws.Range("G1").Value = "soil_texture_sub"

# Fill in the soil_texture_main (general texture group) values for each
# sample; soil_texture_sub already holds the previous, more specific
# soil_texture values (shifted from column F into column G).
$ws.Range("F2").Value = "Clay"
$ws.Range("F3").Value = "Sand"
$ws.Range("F4").Value = "Clay"
$ws.Range("F5").Value = "Sand"

# The new column should render with the same best-fit width the
# soil_texture_sub column (shifted from the old soil_texture column) has.
$ws.Range("F1").ColumnWidth = $ws.Range("G1").ColumnWidth()

# Match the saved selection state seen in the authored workbook.
$ws.Range("F5").Select()
